$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a range to be written/stored as text (matches the source
# workbook's inlineStr cells, which include plain-number-looking strings
# like "1.00" that Excel would otherwise auto-convert to a number).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "59.055.36"
$ws.Range("E2").Value = "  +1.86%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.199.91"
$ws.Range("E3").Value = "  +3.34%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "537.07"
$ws.Range("E5").Value = "  +2.09%  "

# Row 6
Set-TextValue $ws.Range("D6") "144.21"
$ws.Range("E6").Value = "  +2.03%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.194.36"
$ws.Range("E8").Value = "  +3.14%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.452"
$ws.Range("E9").Value = "  +3.15%  "

# Row 10
Set-TextValue $ws.Range("D10") "7.24"
$ws.Range("E10").Value = "  -0.85%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.112"
$ws.Range("E11").Value = "  +2.90%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.402"
$ws.Range("E12").Value = "  +4.75%  "

# Row 13
Set-TextValue $ws.Range("D13") "3.743.99"
$ws.Range("E13").Value = "  +3.25%  "

# Row 14
$ws.Range("E14").Value = "  +3.06%  "

# Row 15
Set-TextValue $ws.Range("D15") "26.23"
$ws.Range("E15").Value = "  -1.51%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.0000170"
$ws.Range("E16").Value = "  +2.78%  "

# Row 17
Set-TextValue $ws.Range("D17") "59.085.38"
$ws.Range("E17").Value = "  +1.78%  "

# Row 18
Set-TextValue $ws.Range("D18") "3.193.17"
$ws.Range("E18").Value = "  +3.19%  "

# Row 19
Set-TextValue $ws.Range("D19") "6.23"
$ws.Range("E19").Value = "  +2.15%  "

# Row 20
Set-TextValue $ws.Range("D20") "13.09"
$ws.Range("E20").Value = "  +1.98%  "

# Row 21
Set-TextValue $ws.Range("D21") "8.15"
$ws.Range("E21").Value = "  +0.97%  "

# Row 22
Set-TextValue $ws.Range("D22") "364.24"
$ws.Range("E22").Value = "  +7.60%  "

# Row 23
$ws.Range("E23").Value = "  -0.11%  "

# Row 24
Set-TextValue $ws.Range("D24") "0.522"
$ws.Range("E24").Value = "  +3.86%  "

# Row 25
Set-TextValue $ws.Range("D25") "68.93"
$ws.Range("E25").Value = "  +4.96%  "

# Row 26
$ws.Range("E26").Value = "  +1.64%  "

# Row 27 (was Binance-PegBSC-USD, now PEPE)
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D27") "0.0₃0968"
$ws.Range("E27").Value = "  +6.73%  "

# Row 28 (was PEPE, now Binance-PegBSC-USD)
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D28") "1.00"
$ws.Range("E28").Value = "  -0.13%  "

# Row 29 (was InternetComputer(DFINITY), now RenderToken)
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D29") "6.63"
$ws.Range("E29").Value = "  +0.95%  "

# Row 30 (was RenderToken, now InternetComputer(DFINITY))
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D30") "7.61"
$ws.Range("E30").Value = "  +4.85%  "

# Row 31
$ws.Range("E31").Value = "  -0.01%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.92"
$ws.Range("E32").Value = "  +3.32%  "

# Row 33
Set-TextValue $ws.Range("D33") "21.65"
$ws.Range("E33").Value = "  +3.78%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.22"
$ws.Range("E34").Value = "  +2.51%  "

# Row 35
Set-TextValue $ws.Range("D35") "4.92"
$ws.Range("E35").Value = "  +6.94%  "

# Row 36
$ws.Range("E36").Value = "  +3.03%  "

# Row 37
Set-TextValue $ws.Range("D37") "6.36"
$ws.Range("E37").Value = "  +5.55%  "

# Row 38
Set-TextValue $ws.Range("D38") "26.95"
$ws.Range("E38").Value = "  -0.33%  "

# Row 39
Set-TextValue $ws.Range("D39") "1.32"
$ws.Range("E39").Value = "  +2.26%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.68"
$ws.Range("E40").Value = "  +13.93%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.0685"
$ws.Range("E41").Value = "  +2.13%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.715"
$ws.Range("E42").Value = "  +5.86%  "

# Row 43
Set-TextValue $ws.Range("D43") "4.08"
$ws.Range("E43").Value = "  +5.51%  "

# Row 44
Set-TextValue $ws.Range("D44") "3.235.27"
$ws.Range("E44").Value = "  +3.08%  "

# Row 45
Set-TextValue $ws.Range("D45") "37.19"
$ws.Range("E45").Value = "  +1.32%  "

# Row 46 (was FirstDigitalUSD, now VeChain)
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D46") "0.0273"
$ws.Range("E46").Value = "  +6.71%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.374.00"
$ws.Range("E47").Value = "  +3.90%  "

# Row 48 (was VeChain, now FirstDigitalUSD)
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D48") "1.00"
$ws.Range("E48").Value = "  +0.10%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.04"
$ws.Range("E49").Value = "  +8.69%  "

# Row 50
Set-TextValue $ws.Range("D50") "21.01"
$ws.Range("E50").Value = "  +1.10%  "

# Row 51
Set-TextValue $ws.Range("D51") "6.14"
$ws.Range("E51").Value = "  +2.93%  "
